# error solve ifrs list
# Recomputed financial figures (매출액/영업이익/... etc.) for 일성건설 rows 2-6,
# and dropped the (E) estimate rows' stale bulk figures down to only the
# identifying columns (A/B/C), matching the corrected IFRS export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 : 2014/12 (IFRS연결) ----
$ws.Range("D2").Value  = 1921
$ws.Range("E2").Value  = -293
$ws.Range("F2").Value  = -293
$ws.Range("G2").Value  = -367
$ws.Range("H2").Value  = -374
$ws.Range("I2").Value  = -374
$ws.Range("J2").Value  = 0
$ws.Range("K2").Value  = 2647
$ws.Range("L2").Value  = 1759
$ws.Range("M2").Value  = 888
$ws.Range("N2").Value  = 889
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 270
$ws.Range("Q2").Value  = -147
$ws.Range("R2").Value  = -69
$ws.Range("S2").Value  = 92
$ws.Range("T2").Value  = 21
$ws.Range("U2").Value  = -167
$ws.Range("V2").Value  = 980
$ws.Range("W2").Value  = -15.25
$ws.Range("X2").Value  = -19.48
$ws.Range("Y2").Value  = -34.58
$ws.Range("Z2").Value  = -13.62
$ws.Range("AA2").Value = 197.99
$ws.Range("AB2").Value = 228.65
$ws.Range("AC2").Value = -693
$ws.Range("AD2").Value = -0.82
$ws.Range("AE2").Value = 1773
$ws.Range("AF2").Value = 0.32
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 53403640

# ---- Row 3 : 2015/12 (IFRS연결) ----
$ws.Range("D3").Value  = 2311
$ws.Range("E3").Value  = 73
$ws.Range("F3").Value  = 69
$ws.Range("G3").Value  = 40
$ws.Range("H3").Value  = 37
$ws.Range("I3").Value  = 37
$ws.Range("J3").Value  = -1
$ws.Range("K3").Value  = 2590
$ws.Range("L3").Value  = 1673
$ws.Range("M3").Value  = 917
$ws.Range("N3").Value  = 917
$ws.Range("O3").Value  = -1
$ws.Range("P3").Value  = 270
$ws.Range("Q3").Value  = -119
$ws.Range("R3").Value  = 185
$ws.Range("S3").Value  = 9
$ws.Range("T3").Value  = 9
$ws.Range("U3").Value  = -128
$ws.Range("V3").Value  = 995
$ws.Range("W3").Value  = 3.17
$ws.Range("X3").Value  = 1.6
$ws.Range("Y3").Value  = 4.15
$ws.Range("Z3").Value  = 1.41
$ws.Range("AA3").Value = 182.57
$ws.Range("AB3").Value = 241.73
$ws.Range("AC3").Value = 69
$ws.Range("AD3").Value = 15.79
$ws.Range("AE3").Value = 1830
$ws.Range("AF3").Value = 0.6
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0.41
$ws.Range("AJ3").Value = 53403640

# ---- Row 4 : 2016/12 (IFRS연결) ----
$ws.Range("D4").Value  = 2858
$ws.Range("E4").Value  = 81
$ws.Range("F4").Value  = 81
$ws.Range("G4").Value  = 23
$ws.Range("H4").Value  = 21
$ws.Range("I4").Value  = 21
$ws.Range("J4").Value  = 0
$ws.Range("K4").Value  = 2775
$ws.Range("L4").Value  = 1800
$ws.Range("M4").Value  = 974
$ws.Range("N4").Value  = 975
$ws.Range("O4").Value  = -1
$ws.Range("P4").Value  = 270
$ws.Range("Q4").Value  = 211
$ws.Range("R4").Value  = -84
$ws.Range("S4").Value  = -123
$ws.Range("T4").Value  = 6
$ws.Range("U4").Value  = 205
$ws.Range("V4").Value  = 865
$ws.Range("W4").Value  = 2.84
$ws.Range("X4").Value  = 0.73
$ws.Range("Y4").Value  = 2.21
$ws.Range("Z4").Value  = 0.77
$ws.Range("AA4").Value = 184.81
$ws.Range("AB4").Value = 253.22
$ws.Range("AC4").Value = 39
$ws.Range("AD4").Value = 22.79
$ws.Range("AE4").Value = 1805
$ws.Range("AF4").Value = 0.49
$ws.Range("AG4").Value = 20
$ws.Range("AH4").Value = 2.27
$ws.Range("AI4").Value = 51.95
$ws.Range("AJ4").Value = 53403640

# ---- Row 5 : 2017/12 (IFRS연결) ----
$ws.Range("D5").Value  = 4208
$ws.Range("E5").Value  = 129
$ws.Range("F5").Value  = 129
$ws.Range("G5").Value  = -46
$ws.Range("H5").Value  = -120
$ws.Range("I5").Value  = -120
$ws.Range("J5").Value  = 0
$ws.Range("K5").Value  = 3053
$ws.Range("L5").Value  = 2211
$ws.Range("M5").Value  = 842
$ws.Range("N5").Value  = 843
$ws.Range("O5").Value  = -1
$ws.Range("P5").Value  = 270
$ws.Range("Q5").Value  = -34
$ws.Range("R5").Value  = -97
$ws.Range("S5").Value  = 119
$ws.Range("T5").Value  = 3
$ws.Range("U5").Value  = -37
$ws.Range("V5").Value  = 992
$ws.Range("W5").Value  = 3.07
$ws.Range("X5").Value  = -2.86
$ws.Range("Y5").Value  = -13.21
$ws.Range("Z5").Value  = -4.13
$ws.Range("AA5").Value = 262.57
$ws.Range("AB5").Value = 205.22
$ws.Range("AC5").Value = -222
$ws.Range("AD5").Value = -4.63
$ws.Range("AE5").Value = 1560
$ws.Range("AF5").Value = 0.66
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = -0.13
$ws.Range("AJ5").Value = 53403640

# ---- Row 6 : 2018/12 (IFRS연결) ----
$ws.Range("D6").Value  = 3640
$ws.Range("E6").Value  = 99
$ws.Range("F6").Value  = 99
$ws.Range("G6").Value  = 30
$ws.Range("H6").Value  = 24
$ws.Range("I6").Value  = 24
$ws.Range("K6").Value  = 3087
$ws.Range("L6").Value  = 2187
$ws.Range("M6").Value  = 901
$ws.Range("N6").Value  = 902
$ws.Range("P6").Value  = 270
$ws.Range("Q6").Value  = 215
$ws.Range("R6").Value  = -108
$ws.Range("S6").Value  = -103
$ws.Range("T6").Value  = 3
$ws.Range("U6").Value  = 212
$ws.Range("V6").Value  = 824
$ws.Range("W6").Value  = 2.71
$ws.Range("X6").Value  = 0.67
$ws.Range("Y6").Value  = 2.79
$ws.Range("Z6").Value  = 0.79
$ws.Range("AA6").Value = 242.71
$ws.Range("AB6").Value = 211.14
$ws.Range("AC6").Value = 45
$ws.Range("AD6").Value = 31.31
$ws.Range("AE6").Value = 1669
$ws.Range("AF6").Value = 0.84
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 54024880

# Row 6 no longer reports 현금DPS(원)/현금배당수익률 (AG6/AH6) at all -
# clear them so the cells disappear from the sheet entirely.
$ws.Range("AG6:AH6").ClearContents()

# Rows 7-9 (the 2019/12(E)-2021/12(E) estimate rows) keep only their
# identifying columns (A row#, B "연간", C period label); every other
# figure (D..AI) is removed.
$ws.Range("D7:AI9").ClearContents()
